$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that are no longer populated in the updated rows
$ws.Range("T2").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("G12").ClearContents()

# Set cell values for rows 2-33 (aging stock summary data)
# Row 2
$ws.Range("A2").Value = '16 - 30 Days'
$ws.Range("B2").Value = 'Osticare'
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 'Osticare Tablet 30''s'
$ws.Range("J2").Value = 1
$ws.Range("Y2").Value = 1

# Row 3
$ws.Range("A3").Value = '16 - 30 Days'
$ws.Range("B3").Value = 'Flucloxin'
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 'Flucloxin 100ml Dry Suspension'
$ws.Range("G3").Value = 22

# Row 4
$ws.Range("A4").Value = '16 - 30 Days'
$ws.Range("B4").Value = 'Mebidal'
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 'Mebidal Tablet'
$ws.Range("AC4").Value = 24

# Row 5
$ws.Range("A5").Value = '16 - 30 Days'
$ws.Range("B5").Value = 'Naprox'
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 'Naprox Plus 500mg Tablet - 36''s'
$ws.Range("Y5").Value = 1

# Row 6
$ws.Range("A6").Value = '16 - 30 Days'
$ws.Range("B6").Value = 'Flucloxin'
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 'Flucloxin 500mg Capsule 40''s'
$ws.Range("G6").Value = 7
$ws.Range("Y6").Value = 1

# Row 7
$ws.Range("A7").Value = '31 - 60 Days'
$ws.Range("B7").Value = 'Oradin'
$ws.Range("C7").Value = 6
$ws.Range("D7").Value = 'Oradin 60ml Suspension'
$ws.Range("J7").Value = 3

# Row 8
$ws.Range("A8").Value = '31 - 60 Days'
$ws.Range("B8").Value = 'Toti'
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 'Toti 100ml Syrup'
$ws.Range("G8").Value = 14

# Row 9
$ws.Range("A9").Value = '61 - 90 Days'
$ws.Range("B9").Value = 'Toti'
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 'Toti 100ml Syrup'
$ws.Range("J9").Value = 5
$ws.Range("AA9").Value = 33

# Row 10
$ws.Range("A10").Value = '61 - 90 Days'
$ws.Range("B10").Value = 'Zithrox'
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 'Zithrox 35ml Dry Suspension'
$ws.Range("J10").Value = 7
$ws.Range("Y10").Value = 1

# Row 11
$ws.Range("A11").Value = '61 - 90 Days'
$ws.Range("B11").Value = 'Dinafex'
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 'Dinafex 50ml Suspension'
$ws.Range("AC11").Value = 1

# Row 12
$ws.Range("A12").Value = '91 - 180 Days'
$ws.Range("B12").Value = 'Dorenta'
$ws.Range("C12").Value = 11
$ws.Range("D12").Value = 'Dorenta 100ml Syrup'
$ws.Range("J12").Value = 42

# Row 13
$ws.Range("A13").Value = '91 - 180 Days'
$ws.Range("B13").Value = 'Zithrox'
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 'Zithrox 35ml Dry Suspension'
$ws.Range("J13").Value = 1

# Row 14
$ws.Range("A14").Value = '91 - 180 Days'
$ws.Range("B14").Value = 'Zithrox'
$ws.Range("C14").Value = 13
$ws.Range("D14").Value = 'Zithrox 20ml Powder for Suspension'
$ws.Range("J14").Value = 3

# Row 15
$ws.Range("A15").Value = '91 - 180 Days'
$ws.Range("B15").Value = 'Flucloxin'
$ws.Range("C15").Value = 14
$ws.Range("D15").Value = 'Flucloxin 100ml Dry Suspension'
$ws.Range("Q15").Value = 4

# Row 16
$ws.Range("A16").Value = '91 - 180 Days'
$ws.Range("B16").Value = 'Mebidal'
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 'Mebidal Tablet'
$ws.Range("Y16").Value = 20

# Row 17
$ws.Range("A17").Value = '91 - 180 Days'
$ws.Range("B17").Value = 'Oradin'
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 'Oradin 60ml Suspension'
$ws.Range("AA17").Value = 1

# Row 18
$ws.Range("A18").Value = '91 - 180 Days'
$ws.Range("B18").Value = 'Toti'
$ws.Range("C18").Value = 17
$ws.Range("D18").Value = 'Toti 100ml Syrup'
$ws.Range("J18").Value = 2
$ws.Range("X18").Value = 5

# Row 19
$ws.Range("A19").Value = '91 - 180 Days'
$ws.Range("B19").Value = 'Augment'
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 'Augment 100ml PFS'
$ws.Range("I19").Value = 1

# Row 20
$ws.Range("A20").Value = '91 - 180 Days'
$ws.Range("B20").Value = 'Naprox'
$ws.Range("C20").Value = 19
$ws.Range("D20").Value = 'Naprox Plus 500mg Tablet - 36''s'
$ws.Range("H20").Value = 78

# Row 21
$ws.Range("A21").Value = '91 - 180 Days'
$ws.Range("B21").Value = 'Ontin'
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 'Ontin 60ml Syrup'
$ws.Range("AC21").Value = 2

# Row 22
$ws.Range("A22").Value = '91 - 180 Days'
$ws.Range("B22").Value = 'Augment'
$ws.Range("C22").Value = 21
$ws.Range("D22").Value = 'Augment 1.2g IV Injection 1''s'
$ws.Range("AD22").Value = 32

# Row 23
$ws.Range("A23").Value = '91 - 180 Days'
$ws.Range("B23").Value = 'Ketonic'
$ws.Range("C23").Value = 22
$ws.Range("D23").Value = 'Ketonic 60mg IM Injection'
$ws.Range("AA23").Value = 6

# Row 24
$ws.Range("A24").Value = '91 - 180 Days'
$ws.Range("B24").Value = 'Ketonic'
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 'Ketonic 30mg IM/IV Injection'
$ws.Range("X24").Value = 1

# Row 25
$ws.Range("A25").Value = '91 - 180 Days'
$ws.Range("B25").Value = 'Flucloxin'
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 'Flucloxin 500mg Capsule 40''s'
$ws.Range("P25").Value = 1
$ws.Range("Q25").Value = 10

# Row 26
$ws.Range("A26").Value = '91 - 180 Days'
$ws.Range("B26").Value = 'Quinox'
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 'Quinox 100ml IV Infusion'
$ws.Range("P26").Value = 4
$ws.Range("R26").Value = 13

# Row 27
$ws.Range("A27").Value = '181 - 210 Days'
$ws.Range("B27").Value = 'Etorix'
$ws.Range("C27").Value = 26
$ws.Range("D27").Value = 'Etorix 60mg Tablet 50''s'
$ws.Range("I27").Value = 1

# Row 28
$ws.Range("A28").Value = '181 - 210 Days'
$ws.Range("B28").Value = 'Quinox'
$ws.Range("C28").Value = 27
$ws.Range("D28").Value = 'Quinox 100ml IV Infusion'
$ws.Range("G28").Value = 23
$ws.Range("O28").Value = 10
$ws.Range("S28").Value = 6
$ws.Range("T28").Value = 6

# Row 29
$ws.Range("A29").Value = '181 - 210 Days'
$ws.Range("B29").Value = 'Etorix'
$ws.Range("C29").Value = 28
$ws.Range("D29").Value = 'Etorix 90mg Tablet 40''s'
$ws.Range("I29").Value = 1

# Row 30
$ws.Range("A30").Value = '181 - 210 Days'
$ws.Range("B30").Value = 'Visomox'
$ws.Range("C30").Value = 29
$ws.Range("D30").Value = 'Visomox 400mg FC Tablet 10''s'
$ws.Range("E30").Value = 79
$ws.Range("G30").Value = 191
$ws.Range("H30").Value = 37
$ws.Range("J30").Value = 167
$ws.Range("L30").Value = 215
$ws.Range("N30").Value = 90
$ws.Range("R30").Value = 192
$ws.Range("U30").Value = 219
$ws.Range("V30").Value = 66
$ws.Range("AF30").Value = 80

# Row 31
$ws.Range("A31").Value = '181 - 210 Days'
$ws.Range("B31").Value = 'Ontin'
$ws.Range("C31").Value = 30
$ws.Range("D31").Value = 'Ontin 10mg Tablet'
$ws.Range("AI31").Value = 1

# Row 32
$ws.Range("A32").Value = '181 - 210 Days'
$ws.Range("B32").Value = 'Oradin'
$ws.Range("C32").Value = 31
$ws.Range("D32").Value = 'Oradin 60ml Suspension'
$ws.Range("Q32").Value = 17

# Row 33
$ws.Range("A33").Value = '181 - 210 Days'
$ws.Range("B33").Value = 'Mebidal'
$ws.Range("C33").Value = 32
$ws.Range("D33").Value = 'Mebidal Tablet'
$ws.Range("P33").Value = 29
$ws.Range("Q33").Value = 5
$ws.Range("S33").Value = 28
$ws.Range("AC33").Value = 35
